{"js": "// Load all paragraphs in the body up front (indices are pinned at load time).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- 1) Update the title paragraph (paragraph 0) ---\n// It contains two runs of text separated by a manual line break:\n//   \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 - 09.02.25\" <br/> \"Why Is Anything Conscious?\"\n// Use search+replace on the whole body so we don't disturb the <w:br/>.\nconst dateSearch = context.document.body.search(\"09.02.25\", { matchCase: true });\ndateSearch.load(\"items\");\nawait context.sync();\nif (dateSearch.items.length > 0) {\n  dateSearch.items[0].insertText(\"08.02.25\", \"Replace\");\n} else {\n  // Fallback: set whole paragraph text if the date string wasn't found verbatim.\n  paragraphs.items[0].insertText(\n    \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 - 08.02.25\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n\nconst titleSearch = context.document.body.search(\"Why Is Anything Conscious?\", {\n  matchCase: true,\n});\ntitleSearch.load(\"items\");\nawait context.sync();\nif (titleSearch.items.length > 0) {\n  titleSearch.items[0].insertText(\n    \"Rejection Sampling IMLE: Designing Priors for Better Few-Shot Image Synthesis\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n\n// --- 2) Replace the text of paragraphs 1-6 ---\nparagraphs.items[1].insertText(\n  \"\u05d4\u05d9\u05d5\u05dd \u05e2\u05d5\u05e9\u05d9\u05dd \u05d4\u05e4\u05e1\u05e7\u05d4 \u05e7\u05dc\u05d4 \u05e2\u05dd LLMs \u05d5\u05e1\u05d5\u05e7\u05e8\u05d9\u05dd \u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05e6\u05d9\u05e2 \u05e9\u05d9\u05d8\u05d4 \u05de\u05e2\u05e0\u05d9\u05d9\u05e0\u05ea \u05dc\u05d0\u05d9\u05de\u05d5\u05df \u05de\u05d5\u05d3\u05dc\u05d9 \u05d2\u05e0\u05e8\u05d8\u05d9\u05d1\u05d9\u05d9\u05dd \u05d1\u05de\u05e7\u05e8\u05d4 \u05e9\u05d9\u05e9 \u05dc\u05db\u05dd \u05de\u05e2\u05d8 \u05d3\u05d0\u05d8\u05d4 \u05dc\u05d0\u05d9\u05de\u05d5\u05df. \u05db\u05d9\u05d3\u05d5\u05e2 \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d2\u05e0\u05e8\u05d8\u05d9\u05d1\u05d9\u05d9\u05dd \u05de\u05d5\u05d3\u05e8\u05e0\u05d9\u05d9\u05dd \u05db\u05de\u05d5 \u05de\u05d5\u05d3\u05dc\u05d9 \u05d3\u05d9\u05e4\u05d5\u05d6\u05d9\u05d4, \u05d2\u05d0\u05e0\u05d9\u05dd, VAEs \u05de\u05e6\u05e8\u05d9\u05db\u05d9\u05dd \u05db\u05de\u05d5\u05ea \u05e2\u05e6\u05d5\u05de\u05d4 \u05e9\u05dc \u05d3\u05d0\u05d8\u05d4 \u05d0\u05d1\u05dc \u05dc\u05e4\u05e2\u05de\u05d9\u05dd \u05d0\u05d9\u05df \u05dc\u05e0\u05d5 \u05d0\u05ea \u05d4\u05dc\u05d5\u05e7\u05e1\u05d5\u05e1 \u05d4\u05d6\u05d4 \u05d5\u05d0\u05e0\u05d5 \u05e6\u05e8\u05d9\u05db\u05d9\u05dd \u05dc\u05d0\u05de\u05df \u05e2\u05dc \u05db\u05de\u05d5\u05ea \u05e7\u05d8\u05e0\u05d4 \u05e9\u05dc \u05d3\u05d0\u05d8\u05d4. \u05d4\u05d0\u05dd \u05d6\u05d4 \u05d0\u05e4\u05e9\u05e8\u05d9 \u05d1\u05db\u05dc\u05dc?\",\n  \"Replace\"\n);\n\nparagraphs.items[2].insertText(\n  \"\u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05e2\u05dc \u05db\u05da \u05d7\u05d9\u05d5\u05d1\u05d9\u05ea (\u05dc\u05e4\u05d7\u05d5\u05ea \u05dc\u05e4\u05d9 \u05d4\u05de\u05d0\u05de\u05e8). \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05e2\u05d9\u05dd \u05e9\u05d9\u05d8\u05d4 \u05d4\u05e0\u05e7\u05e8\u05d0\u05ea RS-IMLE \u05dc\u05d0\u05d9\u05de\u05d5\u05df \u05de\u05d5\u05d3\u05dc \u05d2\u05e0\u05e8\u05d8\u05d9\u05d1\u05d9 \u05e2\u05dd \u05de\u05e2\u05d8 \u05d3\u05d0\u05d8\u05d4 \u05e9\u05de\u05e9\u05db\u05dc\u05dc \u05e9\u05d9\u05d8\u05ea IMLE \u05e9\u05d6\u05d4 Implicit Maximum Likelihood Estimation. \u05d1\u05d2\u05d3\u05d5\u05dc \u05de\u05d0\u05d5\u05d3 IMLE \u05d3\u05d9 \u05d3\u05d5\u05de\u05d4 \u05dc\u05e9\u05d9\u05d8\u05d4 \u05d2\u05e0\u05e8\u05d8\u05d9\u05d1\u05d9\u05ea \u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9\u05ea - \u05d4\u05d9\u05d0 \u05d3\u05d5\u05d2\u05de\u05ea \u05de\u05e9\u05ea\u05e0\u05d4 \u05d1\u05e2\u05dc \u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05e7\u05dc\u05d4 \u05dc\u05d3\u05d2\u05d9\u05de\u05d4 (\u05d2\u05d0\u05d5\u05e1\u05d9\u05ea) z \u05d5\u05de\u05d0\u05de\u05e0\u05ea \u05de\u05d5\u05d3\u05dc \u05d2\u05e0\u05e8\u05d8\u05d9\u05d1\u05d9(\u05e8\u05e9\u05ea \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd) \u05db\u05d3\u05d9 \u05dc\u05d2\u05e0\u05e8\u05d8 \u05e4\u05d9\u05e1\u05ea \u05d3\u05d0\u05d8\u05d4. \u05d4\u05d4\u05d1\u05d3\u05dc \u05d4\u05d5\u05d0 \u05d1\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05dc\u05d5\u05e1: \u05e2\u05dd IMLE \u05dc\u05db\u05dc \u05d3\u05d2\u05d9\u05de\u05d4 x \u05de\u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05d0\u05e0\u05d5 \u05de\u05de\u05d6\u05e2\u05e8\u05d9\u05dd \u05d0\u05ea \u05e8\u05e7 \u05d4\u05de\u05e8\u05d7\u05e7 \u05d1\u05d9\u05e0\u05d4 \u05dc\u05d1\u05d9\u05df \u05e0\u05e7\u05d5\u05d3\u05d4 z_i \u05d0\u05d7\u05ea \u05d1\u05dc\u05d1\u05d3: \u05db\u05d6\u05d5  \u05e9-(T(z_i \u05e9\u05dc\u05d4 \u05d4\u05d9\u05e0\u05d5 \u05e7\u05e8\u05d5\u05d1 \u05d1\u05d9\u05d5\u05ea\u05e8 \u05d0\u05dc\u05d9\u05d4. \u05db\u05d0\u05df (T(z_i \u05d4\u05d9\u05d0 \u05e4\u05d9\u05e1\u05ea \u05d3\u05d0\u05d8\u05d4 \u05e9\u05d2\u05d5\u05e0\u05e8\u05d8\u05d4 \u05de-z_i \u05d5- T \u05d6\u05d4 \u05d4\u05de\u05d5\u05d3\u05dc \u05e9\u05d0\u05e0\u05d5 \u05de\u05d0\u05de\u05e0\u05d9\u05dd.\",\n  \"Replace\"\n);\n\nparagraphs.items[3].insertText(\n  \"\u05db\u05dc\u05d5\u05de\u05e8 \u05d1\u05e9\u05dc\u05d1 \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05e9\u05dc IMLE \u05d0\u05e0\u05d5 \u05d3\u05d5\u05d2\u05de\u05d9\u05dd m \u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05d5\u05de\u05e2\u05d1\u05d9\u05e8\u05d9\u05dd \u05d0\u05d5\u05ea\u05dd \u05d3\u05e8\u05da \u05de\u05d5\u05d3\u05dc T(\u05e0\u05e7\u05e8\u05d0 \u05dc\u05d5 \u05de\u05d9\u05e4\u05d5\u05d9 \u05d1\u05d4\u05de\u05e9\u05da) \u05d5\u05d1\u05d5\u05e0\u05d9\u05dd m \u05e4\u05d9\u05e1\u05d5\u05ea \u05d3\u05d0\u05d8\u05d4 \u05de\u05d2\u05d5\u05e0\u05e8\u05d8\u05d5\u05ea. \u05dc\u05d0\u05d7\u05e8 \u05de\u05db\u05df \u05dc\u05db\u05dc \u05d3\u05d2\u05d9\u05de\u05d4 x_j \u05de\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05d0\u05e0\u05d5 \u05d1\u05d5\u05d7\u05e8\u05d9\u05dd \u05d0\u05ea z_i \u05d4\u05e7\u05e8\u05d5\u05d1\u05d4 \u05d1\u05d9\u05d5\u05ea\u05e8 \u05dc-x_j. \u05d1\u05e1\u05d5\u05e3 \u05e8\u05e7 \u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05db\u05d0\u05dc\u05d5 \u05de\u05e9\u05ea\u05ea\u05e4\u05d5\u05ea \u05d1\u05de\u05d6\u05e2\u05d5\u05e8 \u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05ea \u05dc\u05d5\u05e1. \u05db\u05de\u05d5\u05d1\u05df \u05e9\u05de\u05e1\u05e4\u05e8 \u05d4\u05e0\u05e7\u05d5\u05d3\u05d5\u05ea m \u05d4\u05de\u05d2\u05d5\u05e0\u05e8\u05d8\u05d5\u05ea \u05d1\u05e9\u05dc\u05d1 \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05e6\u05e8\u05d9\u05da \u05dc\u05d4\u05d9\u05d5\u05ea \u05d2\u05d1\u05d5\u05d4 \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05ea \u05de\u05d0\u05e9\u05e8 \u05d2\u05d5\u05d3\u05dc \u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05dc\u05d0\u05d9\u05de\u05d5\u05df n. \u05d4\u05de\u05d8\u05e8\u05d4 \u05e9\u05dc \u05e9\u05d9\u05d8\u05ea \u05d0\u05d9\u05de\u05d5\u05df \u05d6\u05d5 \u05d4\u05d9\u05d0 \u05dc\u05d0\u05e4\u05d8\u05dd \u05d0\u05ea \u05d4\u05de\u05d5\u05d3\u05dc \u05e8\u05e7 \u05e2\u05d1\u05d5\u05e8 \u05d4\u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05d1\u05de\u05e8\u05d7\u05d1 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 (z) \u05e9\u05d4\u05df \u05d4\u05de\u05de\u05d5\u05e4\u05d5\u05ea \u05e7\u05e8\u05d5\u05d1 \u05dc\u05e0\u05e7\u05d5\u05d3\u05ea \u05de\u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8.\",\n  \"Replace\"\n);\n\nparagraphs.items[4].insertText(\n  \"\u05d4\u05d1\u05e2\u05d9\u05d4 \u05e2\u05dd \u05d4\u05d2\u05d9\u05e9\u05d4 \u05d4\u05d6\u05d5 \u05e9\u05d4\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05e9\u05dc \u05d4\u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05f4\u05d4\u05e0\u05d1\u05d7\u05e8\u05d5\u05ea\u05f4 \u05d1\u05de\u05d4\u05dc\u05da \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05db\u05d1\u05e8 \u05dc\u05d0 \u05d2\u05d0\u05d5\u05e1\u05d9\u05ea \u05e9\u05e2\u05dc\u05d5\u05dc \u05dc\u05d9\u05e6\u05d5\u05e8 \u05dc\u05e0\u05d5 \u05d1\u05e2\u05d9\u05d5\u05ea \u05d1\u05d0\u05d9\u05e0\u05e4\u05e8\u05e0\u05e1 \u05db\u05d9 \u05d0\u05e0\u05d5 \u05db\u05df \u05e8\u05d5\u05e6\u05d9\u05dd \u05dc\u05d3\u05d2\u05d5\u05dd \u05d0\u05ea z \u05de\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05d2\u05d0\u05d5\u05e1\u05d9\u05ea. \u05d4\u05de\u05e8\u05d7\u05e7 \u05d1\u05d9\u05df \u05de\u05d9\u05e4\u05d5\u05d9 T \u05e9\u05dc \u05d3\u05d2\u05d9\u05de\u05d4 \u05d2\u05d0\u05d5\u05e1\u05d9\u05ea \u05de\u05e0\u05e7\u05d5\u05d3\u05d4 \u05de\u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05e9\u05d5\u05e0\u05d4 \u05d1\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05de\u05d6\u05d4 \u05e9\u05dc \u05d4\u05d3\u05d2\u05d9\u05de\u05d4 z \u05d4\u05de\u05de\u05d5\u05e4\u05d4 \u05d4\u05db\u05d9 \u05e7\u05e8\u05d5\u05d1 \u05dc\u05e7\u05d5\u05d3\u05d4 \u05d6\u05d5 (\u05d4\u05d0\u05de\u05ea \u05d6\u05d4 \u05d3\u05d9 \u05d1\u05e8\u05d5\u05e8). \u05d3\u05e8\u05da \u05d0\u05d2\u05d1 \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05d5\u05db\u05d9\u05d7 \u05d0\u05ea \u05d4\u05d8\u05e2\u05e0\u05d4 \u05d4\u05d6\u05d5 \u05d5\u05de\u05e6\u05d9\u05e2 \u05e9\u05d9\u05d8\u05d4 \u05dc\u05d4\u05ea\u05d2\u05d1\u05e8 \u05e2\u05dc \u05d6\u05d4. \",\n  \"Replace\"\n);\n\nparagraphs.items[5].insertText(\n  \"\u05d4\u05e9\u05d9\u05d8\u05d4 \u05e9\u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05e0\u05e8\u05d0\u05d9\u05ea \u05de\u05de\u05e9 \u05e4\u05e9\u05d5\u05d8\u05d4 \u05d0\u05da \u05de\u05d1\u05d5\u05e1\u05e1\u05ea \u05e2\u05dc \u05e0\u05d9\u05ea\u05d5\u05d7 \u05de\u05ea\u05de\u05d8\u05d9 \u05d3\u05d9 \u05de\u05e2\u05de\u05d9\u05e7 \u05e9\u05dc \u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05d9\u05d5\u05ea \u05d4\u05de\u05e8\u05d7\u05e7\u05d9\u05dd. \u05d1\u05e9\u05dc\u05d1 \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05e9\u05dc \u05d4\u05d0\u05d9\u05de\u05d5\u05df (\u05d0\u05d7\u05e8\u05d9 \u05d4\u05d3\u05d2\u05d9\u05de\u05d4 \u05de\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05d2\u05d0\u05d5\u05e1\u05d9\u05ea) \u05d1\u05d5\u05d7\u05e8\u05d9\u05dd \u05d0\u05ea z_i \u05db\u05d0\u05e9\u05e8 \u05e0\u05d5\u05e4\u05dc\u05d9\u05dd \u05d1\u05de\u05e8\u05d7\u05e7 \u05d9\u05d5\u05ea\u05e8 \u05d2\u05d3\u05d5\u05dc \u05de\u05d1\u05d5\u05e2 \u05d0\u05e4\u05e1\u05d9\u05dc\u05d5\u05df \u05de\u05db\u05dc \u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05d1\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05d0\u05d7\u05e8\u05d9 \u05d4\u05de\u05d9\u05e4\u05d5\u05d9 (\u05db\u05dc\u05d5\u05de\u05e8 \u05d9\u05e9 \u05dc\u05e0\u05d5 rejection sampling). \u05dc\u05d0\u05d7\u05e8 \u05de\u05db\u05df, \u05d1\u05d3\u05d5\u05de\u05d4 \u05dc-IMLE, \u05dc\u05db\u05dc \u05e0\u05e7\u05d5\u05d3\u05d4 x \u05d1\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05d1\u05d5\u05d7\u05e8\u05d9\u05dd \u05d0\u05ea z \u05e9\u05d4\u05de\u05d9\u05e4\u05d5\u05d9 \u05e9\u05dc\u05d5 \u05e2\u05dd T \u05e0\u05d5\u05e4\u05dc \u05d4\u05db\u05d9 \u05e7\u05e8\u05d5\u05d1 \u05d0\u05dc\u05d9\u05d4 \u05d5\u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05d0\u05ea T \u05dc\u05de\u05d6\u05e2\u05e8 \u05d0\u05ea \u05d4\u05de\u05e8\u05d7\u05e7 \u05d4\u05de\u05de\u05d5\u05e6\u05e2 \u05d1\u05d9\u05df z-s \u05d4\u05e0\u05d1\u05d7\u05e8\u05d9\u05dd \u05dc\u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05d4\u05e2\u05d5\u05d2\u05df \u05e9\u05dc\u05d4\u05dd. \u05d4\u05d9\u05d9\u05e4\u05e8\u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05d4\u05d7\u05e9\u05d5\u05d1\u05d9\u05dd \u05db\u05d0\u05df \u05d6\u05d4 \u05d0\u05e4\u05e1\u05d9\u05dc\u05d5\u05df \u05d5\u05de\u05e1\u05e4\u05e8 \u05e0\u05e7\u05d5\u05d3\u05d5\u05ea z \u05e9\u05e0\u05d3\u05d2\u05de\u05d5\u05ea. \",\n  \"Replace\"\n);\n\nparagraphs.items[6].insertText(\n  \"\u05d0\u05d9\u05e0\u05d8\u05d5\u05d0\u05d9\u05d8\u05d9\u05d1\u05d9\u05ea \u05d6\u05d4 \u05e2\u05d5\u05d1\u05d3 \u05db\u05d9 \u05de\u05dc\u05db\u05ea\u05d7\u05d9\u05dc\u05d4 \u05d0\u05e0\u05d5 \u05d1\u05d5\u05d7\u05e8\u05d9\u05dd \u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05e8\u05d7\u05d5\u05e7\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 (\u05dc\u05d0\u05d7\u05e8 \u05d4\u05de\u05d9\u05e4\u05d5\u05d9) \u05de\u05d4\u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05d1\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05e9\u05de\u05d0\u05e4\u05e9\u05e8 \u05dc\u05e9\u05de\u05d5\u05e8 \u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05e9\u05dc \u05d4\u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05d4\u05e0\u05d1\u05d7\u05e8\u05d5\u05ea \u05d1\u05e9\u05dc\u05d1 \u05dc\u05d0\u05d7\u05e8 \u05de\u05db\u05df \u05e7\u05e8\u05d5\u05d1\u05d4 \u05dc\u05d2\u05d0\u05d5\u05e1\u05d9\u05ea.\",\n  \"Replace\"\n);\n\nawait context.sync();\n\n// --- 3) Delete paragraphs 7..33 (inclusive) - all of the removed content. ---\n// These index-pinned proxies refer to the paragraphs as loaded at the top,\n// so deleting them in order (front to back) is safe and each call still\n// targets the originally-intended paragraph.\nfor (let i = 7; i <= 33; i++) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n\n// --- 4) Fix up the final paragraph's URL. ---\nconst urlSearch = context.document.body.search(\"2409.14545\", { matchCase: true });\nurlSearch.load(\"items\");\nawait context.sync();\nif (urlSearch.items.length > 0) {\n  urlSearch.items[0].insertText(\"2409.17439\", \"Replace\");\n} else {\n  paragraphs.items[34].insertText(\n    \"https://arxiv.org/abs/2409.17439\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1) Update the title paragraph (paragraph 1) ---\n# It holds two runs of text separated by a manual line break:\n#   \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 - 09.02.25\" <br/> \"Why Is Anything Conscious?\"\n# Use Find/Replace so the <w:br/> between the two runs is left untouched.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"09.02.25\", $false, $false, $false, $false, $false, $true, 1, $false, \"08.02.25\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"Why Is Anything Conscious?\", $false, $false, $false, $false, $false, $true, 1, $false, \"Rejection Sampling IMLE: Designing Priors for Better Few-Shot Image Synthesis\", 2) | Out-Null\n\n# --- 2) Replace the text of paragraphs 2-7 (1-indexed) in place ---\n$d.Paragraphs(2).Range.Text = \"\u05d4\u05d9\u05d5\u05dd \u05e2\u05d5\u05e9\u05d9\u05dd \u05d4\u05e4\u05e1\u05e7\u05d4 \u05e7\u05dc\u05d4 \u05e2\u05dd LLMs \u05d5\u05e1\u05d5\u05e7\u05e8\u05d9\u05dd \u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05e6\u05d9\u05e2 \u05e9\u05d9\u05d8\u05d4 \u05de\u05e2\u05e0\u05d9\u05d9\u05e0\u05ea \u05dc\u05d0\u05d9\u05de\u05d5\u05df \u05de\u05d5\u05d3\u05dc\u05d9 \u05d2\u05e0\u05e8\u05d8\u05d9\u05d1\u05d9\u05d9\u05dd \u05d1\u05de\u05e7\u05e8\u05d4 \u05e9\u05d9\u05e9 \u05dc\u05db\u05dd \u05de\u05e2\u05d8 \u05d3\u05d0\u05d8\u05d4 \u05dc\u05d0\u05d9\u05de\u05d5\u05df. \u05db\u05d9\u05d3\u05d5\u05e2 \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d2\u05e0\u05e8\u05d8\u05d9\u05d1\u05d9\u05d9\u05dd \u05de\u05d5\u05d3\u05e8\u05e0\u05d9\u05d9\u05dd \u05db\u05de\u05d5 \u05de\u05d5\u05d3\u05dc\u05d9 \u05d3\u05d9\u05e4\u05d5\u05d6\u05d9\u05d4, \u05d2\u05d0\u05e0\u05d9\u05dd, VAEs \u05de\u05e6\u05e8\u05d9\u05db\u05d9\u05dd \u05db\u05de\u05d5\u05ea \u05e2\u05e6\u05d5\u05de\u05d4 \u05e9\u05dc \u05d3\u05d0\u05d8\u05d4 \u05d0\u05d1\u05dc \u05dc\u05e4\u05e2\u05de\u05d9\u05dd \u05d0\u05d9\u05df \u05dc\u05e0\u05d5 \u05d0\u05ea \u05d4\u05dc\u05d5\u05e7\u05e1\u05d5\u05e1 \u05d4\u05d6\u05d4 \u05d5\u05d0\u05e0\u05d5 \u05e6\u05e8\u05d9\u05db\u05d9\u05dd \u05dc\u05d0\u05de\u05df \u05e2\u05dc \u05db\u05de\u05d5\u05ea \u05e7\u05d8\u05e0\u05d4 \u05e9\u05dc \u05d3\u05d0\u05d8\u05d4. \u05d4\u05d0\u05dd \u05d6\u05d4 \u05d0\u05e4\u05e9\u05e8\u05d9 \u05d1\u05db\u05dc\u05dc?\"\n\n$d.Paragraphs(3).Range.Text = \"\u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05e2\u05dc \u05db\u05da \u05d7\u05d9\u05d5\u05d1\u05d9\u05ea (\u05dc\u05e4\u05d7\u05d5\u05ea \u05dc\u05e4\u05d9 \u05d4\u05de\u05d0\u05de\u05e8). \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05e2\u05d9\u05dd \u05e9\u05d9\u05d8\u05d4 \u05d4\u05e0\u05e7\u05e8\u05d0\u05ea RS-IMLE \u05dc\u05d0\u05d9\u05de\u05d5\u05df \u05de\u05d5\u05d3\u05dc \u05d2\u05e0\u05e8\u05d8\u05d9\u05d1\u05d9 \u05e2\u05dd \u05de\u05e2\u05d8 \u05d3\u05d0\u05d8\u05d4 \u05e9\u05de\u05e9\u05db\u05dc\u05dc \u05e9\u05d9\u05d8\u05ea IMLE \u05e9\u05d6\u05d4 Implicit Maximum Likelihood Estimation. \u05d1\u05d2\u05d3\u05d5\u05dc \u05de\u05d0\u05d5\u05d3 IMLE \u05d3\u05d9 \u05d3\u05d5\u05de\u05d4 \u05dc\u05e9\u05d9\u05d8\u05d4 \u05d2\u05e0\u05e8\u05d8\u05d9\u05d1\u05d9\u05ea \u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9\u05ea - \u05d4\u05d9\u05d0 \u05d3\u05d5\u05d2\u05de\u05ea \u05de\u05e9\u05ea\u05e0\u05d4 \u05d1\u05e2\u05dc \u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05e7\u05dc\u05d4 \u05dc\u05d3\u05d2\u05d9\u05de\u05d4 (\u05d2\u05d0\u05d5\u05e1\u05d9\u05ea) z \u05d5\u05de\u05d0\u05de\u05e0\u05ea \u05de\u05d5\u05d3\u05dc \u05d2\u05e0\u05e8\u05d8\u05d9\u05d1\u05d9(\u05e8\u05e9\u05ea \u05e0\u05d5\u05d9\u05e8\u05d5\u05e0\u05d9\u05dd) \u05db\u05d3\u05d9 \u05dc\u05d2\u05e0\u05e8\u05d8 \u05e4\u05d9\u05e1\u05ea \u05d3\u05d0\u05d8\u05d4. \u05d4\u05d4\u05d1\u05d3\u05dc \u05d4\u05d5\u05d0 \u05d1\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05dc\u05d5\u05e1: \u05e2\u05dd IMLE \u05dc\u05db\u05dc \u05d3\u05d2\u05d9\u05de\u05d4 x \u05de\u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05d0\u05e0\u05d5 \u05de\u05de\u05d6\u05e2\u05e8\u05d9\u05dd \u05d0\u05ea \u05e8\u05e7 \u05d4\u05de\u05e8\u05d7\u05e7 \u05d1\u05d9\u05e0\u05d4 \u05dc\u05d1\u05d9\u05df \u05e0\u05e7\u05d5\u05d3\u05d4 z_i \u05d0\u05d7\u05ea \u05d1\u05dc\u05d1\u05d3: \u05db\u05d6\u05d5  \u05e9-(T(z_i \u05e9\u05dc\u05d4 \u05d4\u05d9\u05e0\u05d5 \u05e7\u05e8\u05d5\u05d1 \u05d1\u05d9\u05d5\u05ea\u05e8 \u05d0\u05dc\u05d9\u05d4. \u05db\u05d0\u05df (T(z_i \u05d4\u05d9\u05d0 \u05e4\u05d9\u05e1\u05ea \u05d3\u05d0\u05d8\u05d4 \u05e9\u05d2\u05d5\u05e0\u05e8\u05d8\u05d4 \u05de-z_i \u05d5- T \u05d6\u05d4 \u05d4\u05de\u05d5\u05d3\u05dc \u05e9\u05d0\u05e0\u05d5 \u05de\u05d0\u05de\u05e0\u05d9\u05dd.\"\n\n$d.Paragraphs(4).Range.Text = \"\u05db\u05dc\u05d5\u05de\u05e8 \u05d1\u05e9\u05dc\u05d1 \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05e9\u05dc IMLE \u05d0\u05e0\u05d5 \u05d3\u05d5\u05d2\u05de\u05d9\u05dd m \u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05d5\u05de\u05e2\u05d1\u05d9\u05e8\u05d9\u05dd \u05d0\u05d5\u05ea\u05dd \u05d3\u05e8\u05da \u05de\u05d5\u05d3\u05dc T(\u05e0\u05e7\u05e8\u05d0 \u05dc\u05d5 \u05de\u05d9\u05e4\u05d5\u05d9 \u05d1\u05d4\u05de\u05e9\u05da) \u05d5\u05d1\u05d5\u05e0\u05d9\u05dd m \u05e4\u05d9\u05e1\u05d5\u05ea \u05d3\u05d0\u05d8\u05d4 \u05de\u05d2\u05d5\u05e0\u05e8\u05d8\u05d5\u05ea. \u05dc\u05d0\u05d7\u05e8 \u05de\u05db\u05df \u05dc\u05db\u05dc \u05d3\u05d2\u05d9\u05de\u05d4 x_j \u05de\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05d0\u05e0\u05d5 \u05d1\u05d5\u05d7\u05e8\u05d9\u05dd \u05d0\u05ea z_i \u05d4\u05e7\u05e8\u05d5\u05d1\u05d4 \u05d1\u05d9\u05d5\u05ea\u05e8 \u05dc-x_j. \u05d1\u05e1\u05d5\u05e3 \u05e8\u05e7 \u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05db\u05d0\u05dc\u05d5 \u05de\u05e9\u05ea\u05ea\u05e4\u05d5\u05ea \u05d1\u05de\u05d6\u05e2\u05d5\u05e8 \u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05ea \u05dc\u05d5\u05e1. \u05db\u05de\u05d5\u05d1\u05df \u05e9\u05de\u05e1\u05e4\u05e8 \u05d4\u05e0\u05e7\u05d5\u05d3\u05d5\u05ea m \u05d4\u05de\u05d2\u05d5\u05e0\u05e8\u05d8\u05d5\u05ea \u05d1\u05e9\u05dc\u05d1 \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05e6\u05e8\u05d9\u05da \u05dc\u05d4\u05d9\u05d5\u05ea \u05d2\u05d1\u05d5\u05d4 \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05ea \u05de\u05d0\u05e9\u05e8 \u05d2\u05d5\u05d3\u05dc \u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05dc\u05d0\u05d9\u05de\u05d5\u05df n. \u05d4\u05de\u05d8\u05e8\u05d4 \u05e9\u05dc \u05e9\u05d9\u05d8\u05ea \u05d0\u05d9\u05de\u05d5\u05df \u05d6\u05d5 \u05d4\u05d9\u05d0 \u05dc\u05d0\u05e4\u05d8\u05dd \u05d0\u05ea \u05d4\u05de\u05d5\u05d3\u05dc \u05e8\u05e7 \u05e2\u05d1\u05d5\u05e8 \u05d4\u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05d1\u05de\u05e8\u05d7\u05d1 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 (z) \u05e9\u05d4\u05df \u05d4\u05de\u05de\u05d5\u05e4\u05d5\u05ea \u05e7\u05e8\u05d5\u05d1 \u05dc\u05e0\u05e7\u05d5\u05d3\u05ea \u05de\u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8.\"\n\n$d.Paragraphs(5).Range.Text = \"\u05d4\u05d1\u05e2\u05d9\u05d4 \u05e2\u05dd \u05d4\u05d2\u05d9\u05e9\u05d4 \u05d4\u05d6\u05d5 \u05e9\u05d4\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05e9\u05dc \u05d4\u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05f4\u05d4\u05e0\u05d1\u05d7\u05e8\u05d5\u05ea\u05f4 \u05d1\u05de\u05d4\u05dc\u05da \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05db\u05d1\u05e8 \u05dc\u05d0 \u05d2\u05d0\u05d5\u05e1\u05d9\u05ea \u05e9\u05e2\u05dc\u05d5\u05dc \u05dc\u05d9\u05e6\u05d5\u05e8 \u05dc\u05e0\u05d5 \u05d1\u05e2\u05d9\u05d5\u05ea \u05d1\u05d0\u05d9\u05e0\u05e4\u05e8\u05e0\u05e1 \u05db\u05d9 \u05d0\u05e0\u05d5 \u05db\u05df \u05e8\u05d5\u05e6\u05d9\u05dd \u05dc\u05d3\u05d2\u05d5\u05dd \u05d0\u05ea z \u05de\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05d2\u05d0\u05d5\u05e1\u05d9\u05ea. \u05d4\u05de\u05e8\u05d7\u05e7 \u05d1\u05d9\u05df \u05de\u05d9\u05e4\u05d5\u05d9 T \u05e9\u05dc \u05d3\u05d2\u05d9\u05de\u05d4 \u05d2\u05d0\u05d5\u05e1\u05d9\u05ea \u05de\u05e0\u05e7\u05d5\u05d3\u05d4 \u05de\u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05e9\u05d5\u05e0\u05d4 \u05d1\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05de\u05d6\u05d4 \u05e9\u05dc \u05d4\u05d3\u05d2\u05d9\u05de\u05d4 z \u05d4\u05de\u05de\u05d5\u05e4\u05d4 \u05d4\u05db\u05d9 \u05e7\u05e8\u05d5\u05d1 \u05dc\u05e7\u05d5\u05d3\u05d4 \u05d6\u05d5 (\u05d4\u05d0\u05de\u05ea \u05d6\u05d4 \u05d3\u05d9 \u05d1\u05e8\u05d5\u05e8). \u05d3\u05e8\u05da \u05d0\u05d2\u05d1 \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05d5\u05db\u05d9\u05d7 \u05d0\u05ea \u05d4\u05d8\u05e2\u05e0\u05d4 \u05d4\u05d6\u05d5 \u05d5\u05de\u05e6\u05d9\u05e2 \u05e9\u05d9\u05d8\u05d4 \u05dc\u05d4\u05ea\u05d2\u05d1\u05e8 \u05e2\u05dc \u05d6\u05d4. \"\n\n$d.Paragraphs(6).Range.Text = \"\u05d4\u05e9\u05d9\u05d8\u05d4 \u05e9\u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05e0\u05e8\u05d0\u05d9\u05ea \u05de\u05de\u05e9 \u05e4\u05e9\u05d5\u05d8\u05d4 \u05d0\u05da \u05de\u05d1\u05d5\u05e1\u05e1\u05ea \u05e2\u05dc \u05e0\u05d9\u05ea\u05d5\u05d7 \u05de\u05ea\u05de\u05d8\u05d9 \u05d3\u05d9 \u05de\u05e2\u05de\u05d9\u05e7 \u05e9\u05dc \u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05d9\u05d5\u05ea \u05d4\u05de\u05e8\u05d7\u05e7\u05d9\u05dd. \u05d1\u05e9\u05dc\u05d1 \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05e9\u05dc \u05d4\u05d0\u05d9\u05de\u05d5\u05df (\u05d0\u05d7\u05e8\u05d9 \u05d4\u05d3\u05d2\u05d9\u05de\u05d4 \u05de\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05d2\u05d0\u05d5\u05e1\u05d9\u05ea) \u05d1\u05d5\u05d7\u05e8\u05d9\u05dd \u05d0\u05ea z_i \u05db\u05d0\u05e9\u05e8 \u05e0\u05d5\u05e4\u05dc\u05d9\u05dd \u05d1\u05de\u05e8\u05d7\u05e7 \u05d9\u05d5\u05ea\u05e8 \u05d2\u05d3\u05d5\u05dc \u05de\u05d1\u05d5\u05e2 \u05d0\u05e4\u05e1\u05d9\u05dc\u05d5\u05df \u05de\u05db\u05dc \u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05d1\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05d0\u05d7\u05e8\u05d9 \u05d4\u05de\u05d9\u05e4\u05d5\u05d9 (\u05db\u05dc\u05d5\u05de\u05e8 \u05d9\u05e9 \u05dc\u05e0\u05d5 rejection sampling). \u05dc\u05d0\u05d7\u05e8 \u05de\u05db\u05df, \u05d1\u05d3\u05d5\u05de\u05d4 \u05dc-IMLE, \u05dc\u05db\u05dc \u05e0\u05e7\u05d5\u05d3\u05d4 x \u05d1\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05d1\u05d5\u05d7\u05e8\u05d9\u05dd \u05d0\u05ea z \u05e9\u05d4\u05de\u05d9\u05e4\u05d5\u05d9 \u05e9\u05dc\u05d5 \u05e2\u05dd T \u05e0\u05d5\u05e4\u05dc \u05d4\u05db\u05d9 \u05e7\u05e8\u05d5\u05d1 \u05d0\u05dc\u05d9\u05d4 \u05d5\u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05d0\u05ea T \u05dc\u05de\u05d6\u05e2\u05e8 \u05d0\u05ea \u05d4\u05de\u05e8\u05d7\u05e7 \u05d4\u05de\u05de\u05d5\u05e6\u05e2 \u05d1\u05d9\u05df z-s \u05d4\u05e0\u05d1\u05d7\u05e8\u05d9\u05dd \u05dc\u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05d4\u05e2\u05d5\u05d2\u05df \u05e9\u05dc\u05d4\u05dd. \u05d4\u05d9\u05d9\u05e4\u05e8\u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05d4\u05d7\u05e9\u05d5\u05d1\u05d9\u05dd \u05db\u05d0\u05df \u05d6\u05d4 \u05d0\u05e4\u05e1\u05d9\u05dc\u05d5\u05df \u05d5\u05de\u05e1\u05e4\u05e8 \u05e0\u05e7\u05d5\u05d3\u05d5\u05ea z \u05e9\u05e0\u05d3\u05d2\u05de\u05d5\u05ea. \"\n\n$d.Paragraphs(7).Range.Text = \"\u05d0\u05d9\u05e0\u05d8\u05d5\u05d0\u05d9\u05d8\u05d9\u05d1\u05d9\u05ea \u05d6\u05d4 \u05e2\u05d5\u05d1\u05d3 \u05db\u05d9 \u05de\u05dc\u05db\u05ea\u05d7\u05d9\u05dc\u05d4 \u05d0\u05e0\u05d5 \u05d1\u05d5\u05d7\u05e8\u05d9\u05dd \u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05e8\u05d7\u05d5\u05e7\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 (\u05dc\u05d0\u05d7\u05e8 \u05d4\u05de\u05d9\u05e4\u05d5\u05d9) \u05de\u05d4\u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05d1\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05e9\u05de\u05d0\u05e4\u05e9\u05e8 \u05dc\u05e9\u05de\u05d5\u05e8 \u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05e9\u05dc \u05d4\u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05d4\u05e0\u05d1\u05d7\u05e8\u05d5\u05ea \u05d1\u05e9\u05dc\u05d1 \u05dc\u05d0\u05d7\u05e8 \u05de\u05db\u05df \u05e7\u05e8\u05d5\u05d1\u05d4 \u05dc\u05d2\u05d0\u05d5\u05e1\u05d9\u05ea.\"\n\n# --- 3) Delete paragraphs 8..34 (1-indexed) - all of the removed content. ---\n$startRange = $d.Paragraphs(8).Range\n$endRange = $d.Paragraphs(34).Range\n$delRange = $d.Range($startRange.Start, $endRange.End)\n$delRange.Delete()\n\n# --- 4) Fix up the final paragraph's URL. ---\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"2409.14545\", $false, $false, $false, $false, $false, $true, 1, $false, \"2409.17439\", 2) | Out-Null\n"}
